$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new blank row above the data (current row 2), pushing the
#    roster down by one row. Excel's row-insert carries formatting down from
#    the row above for any column that already has a styled cell there.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).Insert()

# The border on A2 doesn't survive the implicit copy-down for the left-most
# column; re-copy A1's format (bold, centered, bordered) onto A2 explicitly.
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Columns H1/H2 had no prior cell content, so nothing got copied down into
# them automatically - stamp them with the same header style used by the
# other header cells (B1/G1/J1) before writing their text.
$ws.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Widen column H (it already existed, just empty) to match column I's
#    10.83-character width (same "characters" setting already used there).
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 10

# ---------------------------------------------------------------------------
# 3) Header row text (order matters for shared-string table insertion order:
#    Prob State, " Approved", SRS, then the REVIEWED BY... banner).
# ---------------------------------------------------------------------------
$ws.Range("H1").Value = "Prob State"
$ws.Range("H2").Value = " Approved"
$ws.Range("J2").Value = "SRS"
$ws.Range("J1").Value = "REVIEWED BY THE PERSON LISTED BELOW"

# ---------------------------------------------------------------------------
# 4) "Prob State" column values per person (column H, rows 3-14).
# ---------------------------------------------------------------------------
$ws.Range("H3").Value = "Yes"
$ws.Range("H4").Value = "Yes"
$ws.Range("H5").Value = "Yes"
$ws.Range("H7").Value = "Yes"
$ws.Range("H8").Value = "Yes"
$ws.Range("H9").Value = "Yes"
$ws.Range("H10").Value = "Yes"
$ws.Range("H11").Value = "Yes"
$ws.Range("H12").Value = "Yes"
$ws.Range("H13").Value = "Yes"
$ws.Range("H14").Value = "Yes"

# ---------------------------------------------------------------------------
# 5) Column J now holds formulas that each point at the *next* row's name
#    (wrapping from the last row back to the first).
# ---------------------------------------------------------------------------
$ws.Range("J3").Formula = "=A4"
$ws.Range("J4").Formula = "=A5"
$ws.Range("J5").Formula = "=A6"
$ws.Range("J6").Formula = "=A7"
$ws.Range("J7").Formula = "=A8"
$ws.Range("J8").Formula = "=A9"
$ws.Range("J9").Formula = "=A10"
$ws.Range("J10").Formula = "=A11"
$ws.Range("J11").Formula = "=A12"
$ws.Range("J12").Formula = "=A13"
$ws.Range("J13").Formula = "=A14"
$ws.Range("J14").Formula = "=A3"

# ---------------------------------------------------------------------------
# 6) Move the active selection to the new "next empty row" cell, J15.
# ---------------------------------------------------------------------------
$ws.Range("J15").Select()
